$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

$ALC.Range("H111").Value = 2013.3334
$ALC.Range("I111").Value = 1300
$ALC.Range("J111").Value = 2726.6667
$ALC.Range("K111").Value = 3900
$ALC.Range("L111").Value = 8180.000100000001
$ALC.Range("M111").Value = -833
$ALC.Range("N111").Value = -14314.0001
$ALC.Range("H133").Value = 43510
$ALC.Range("J133").Value = 43510
$ALC.Range("L133").Value = 43510
$ALC.Range("N133").Value = -53630
$ALC.Range("H137").Value = 40001040
$ALC.Range("I137").Value = 50000890
$ALC.Range("J137").Value = 1646
$ALC.Range("K137").Value = 150002670
$ALC.Range("L137").Value = 4938
$ALC.Range("M137").Value = -150000120
$ALC.Range("N137").Value = -10038
$ARM.Range("H2").Value = 57230.223
$ARM.Range("I2").Value = 78603.38
$ARM.Range("J2").Value = 1660
$ARM.Range("K2").Value = 78603.38
$ARM.Range("L2").Value = 1660
$ARM.Range("M2").Value = -78490.38
$ARM.Range("N2").Value = -1886
$ARM.Range("H32").Value = 23986.266
$ARM.Range("I32").Value = 3006.8044
$ARM.Range("J32").Value = 345671.34
$ARM.Range("K32").Value = 3006.8044
$ARM.Range("L32").Value = 345671.34
$ARM.Range("M32").Value = -2719.8044
$ARM.Range("N32").Value = -346245.34
$ARM.Range("H61").Value = 2573
$ARM.Range("I61").Value = 2045.5714
$ARM.Range("J61").Value = 5210.143
$ARM.Range("K61").Value = 2045.5714
$ARM.Range("L61").Value = 5210.143
$ARM.Range("M61").Value = -1833.5714
$ARM.Range("N61").Value = -5634.143
$ARM.Range("H98").Value = 0
$ARM.Range("J98").Value = 0
$ARM.Range("L98").Value = 0
$ARM.Range("N98").ClearContents()
$ARM.Range("H116").Value = 57230.223
$ARM.Range("I116").Value = 78603.38
$ARM.Range("J116").Value = 1660
$ARM.Range("K116").Value = 78603.38
$ARM.Range("L116").Value = 1660
$ARM.Range("M116").Value = -76309.38
$ARM.Range("N116").Value = -6248
$ARM.Range("H124").Value = 50077.25
$ARM.Range("J124").Value = 50077.25
$ARM.Range("L124").Value = 50077.25
$ARM.Range("N124").Value = -59897.25
$ARM.Range("H125").Value = 30221.5
$ARM.Range("J125").Value = 30246.111
$ARM.Range("L125").Value = 30246.111
$ARM.Range("N125").Value = -40086.111
$ARM.Range("H132").Value = 2467.7222
$ARM.Range("I132").Value = 2084.8696
$ARM.Range("J132").Value = 4669.125
$ARM.Range("K132").Value = 6254.6088
$ARM.Range("L132").Value = 14007.375
$ARM.Range("M132").Value = -3724.6088
$ARM.Range("N132").Value = -19067.375
$ARM.Range("H135").Value = 39090.668
$ARM.Range("J135").Value = 39090.668
$ARM.Range("L135").Value = 39090.668
$ARM.Range("N135").Value = -49230.668
$ARM.Range("H136").Value = 2573
$ARM.Range("I136").Value = 2045.5714
$ARM.Range("J136").Value = 5210.143
$ARM.Range("K136").Value = 6136.7142
$ARM.Range("L136").Value = 15630.429
$ARM.Range("M136").Value = -3586.7142
$ARM.Range("N136").Value = -20730.429
$BSM.Range("H3").Value = 57230.223
$BSM.Range("I3").Value = 78603.38
$BSM.Range("J3").Value = 1660
$BSM.Range("K3").Value = 78603.38
$BSM.Range("L3").Value = 1660
$BSM.Range("M3").Value = -78489.38
$BSM.Range("N3").Value = -1888
$BSM.Range("H59").Value = 0
$BSM.Range("J59").Value = 0
$BSM.Range("L59").Value = 0
$BSM.Range("N59").ClearContents()
$BSM.Range("H92").Value = 0
$BSM.Range("J92").Value = 0
$BSM.Range("L92").Value = 0
$BSM.Range("N92").ClearContents()
$BSM.Range("H100").Value = 33666.668
$BSM.Range("J100").Value = 33666.668
$BSM.Range("L100").Value = 33666.668
$BSM.Range("N100").Value = -35830.668
$BSM.Range("H107").Value = 2311
$BSM.Range("I107").Value = 2234.4443
$BSM.Range("J107").Value = 3000
$BSM.Range("K107").Value = 2234.4443
$BSM.Range("L107").Value = 3000
$BSM.Range("M107").Value = -314.4443000000001
$BSM.Range("N107").Value = -6840
$BSM.Range("H135").Value = 50074.75
$BSM.Range("J135").Value = 50074.75
$BSM.Range("L135").Value = 50074.75
$BSM.Range("N135").Value = -60214.75
$CRP.Range("H92").Value = 29400
$CRP.Range("J92").Value = 29400
$CRP.Range("L92").Value = 29400
$CRP.Range("N92").Value = -34392
$CRP.Range("H107").Value = 575.7727
$CRP.Range("J107").Value = 781.6667
$CRP.Range("L107").Value = 781.6667
$CRP.Range("N107").Value = -4621.6667
$CRP.Range("H132").Value = 2691.6099
$CRP.Range("I132").Value = 2375.4443
$CRP.Range("K132").Value = 7126.3329
$CRP.Range("M132").Value = -4596.3329
$CUL.Range("H35").Value = 3667.3333
$CUL.Range("I35").Value = 1002
$CUL.Range("J35").Value = 5000
$CUL.Range("K35").Value = 3006
$CUL.Range("L35").Value = 15000
$CUL.Range("M35").Value = -2718
$CUL.Range("N35").Value = -15576
$CUL.Range("H131").Value = 6062001
$CUL.Range("I131").Value = 0
$CUL.Range("J131").Value = 6062001
$CUL.Range("K131").Value = 0
$CUL.Range("L131").Value = 18186003
$CUL.Range("N131").Value = -18196083
$CUL.Range("M131").ClearContents()
$GSM.Range("H95").Value = 14875
$GSM.Range("J95").Value = 14875
$GSM.Range("L95").Value = 14875
$GSM.Range("N95").Value = -20367
$GSM.Range("H98").Value = 30643
$GSM.Range("J98").Value = 30643
$GSM.Range("L98").Value = 30643
$GSM.Range("N98").Value = -36633
$GSM.Range("H105").Value = 50671
$GSM.Range("J105").Value = 50671
$GSM.Range("L105").Value = 50671
$GSM.Range("N105").Value = -57659
$GSM.Range("H107").Value = 1277.7778
$GSM.Range("I107").Value = 1899.4
$GSM.Range("K107").Value = 1899.4
$GSM.Range("M107").Value = 20.59999999999991
$GSM.Range("H138").Value = 64147.46
$GSM.Range("J138").Value = 64147.46
$GSM.Range("L138").Value = 64147.46
$GSM.Range("N138").Value = -74427.45999999999
$LTW.Range("H16").Value = 20002590
$LTW.Range("I16").Value = 100000000
$LTW.Range("J16").Value = 3237.5
$LTW.Range("K16").Value = 100000000
$LTW.Range("L16").Value = 3237.5
$LTW.Range("M16").Value = -99999830
$LTW.Range("N16").Value = -3577.5
$LTW.Range("H106").Value = 24148.166
$LTW.Range("J106").Value = 24148.166
$LTW.Range("L106").Value = 24148.166
$LTW.Range("N106").Value = -26672.166
$LTW.Range("H127").Value = 44384.25
$LTW.Range("J127").Value = 44384.25
$LTW.Range("L127").Value = 44384.25
$LTW.Range("N127").Value = -54304.25
$LTW.Range("H136").Value = 5022.25
$LTW.Range("I136").Value = 2821.6765
$LTW.Range("J136").Value = 12504.2
$LTW.Range("K136").Value = 8465.029500000001
$LTW.Range("L136").Value = 37512.60000000001
$LTW.Range("M136").Value = -5915.029500000001
$LTW.Range("N136").Value = -42612.60000000001
$WVR.Range("H69").Value = 15914
$WVR.Range("I69").Value = 0
$WVR.Range("J69").Value = 15914
$WVR.Range("K69").Value = 0
$WVR.Range("L69").Value = 15914
$WVR.Range("M69").ClearContents()
$WVR.Range("N69").Value = -17412
$WVR.Range("H72").Value = 15914
$WVR.Range("I72").Value = 0
$WVR.Range("J72").Value = 15914
$WVR.Range("K72").Value = 0
$WVR.Range("L72").Value = 47742
$WVR.Range("M72").ClearContents()
$WVR.Range("N72").Value = -55230
$WVR.Range("H98").Value = 28663.334
$WVR.Range("J98").Value = 28663.334
$WVR.Range("L98").Value = 28663.334
$WVR.Range("N98").Value = -34653.334
$WVR.Range("H103").Value = 191000.67
$WVR.Range("J103").Value = 191000.67
$WVR.Range("L103").Value = 191000.67
$WVR.Range("N103").Value = -193344.67
